$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheets 1-4: "Fonte/Tecnologia" tables (Potencia Acumulada, Geracao Periodo
# Medio, Atendimento a Ponta, Potencia Incremental).
# Add a header in A1 ("Fonte/Tecnologia"), copying the style already used by
# B1:E1, fix accented labels in column A and drop the (now redundant) header
# style from the row-label cells A2:A12.
# ---------------------------------------------------------------------------
for ($i = 1; $i -le 4; $i++) {
    $ws = $wb.Worksheets.Item($i)

    # Give A1 the same style as the other header cells, then set its text.
    $ws.Range("B1").Copy()
    $ws.Range("A1").PasteSpecial(-4122)
    $excel.CutCopyMode = $false
    $ws.Range("A1").Value = "Fonte/Tecnologia"

    # Row labels: strip the bold header style and fix accents.
    $ws.Range("A2").ClearFormats()
    $ws.Range("A2").Value = "Hidro"

    $ws.Range("A3").ClearFormats()
    $ws.Range("A3").Value = "Gás Natural"

    $ws.Range("A4").ClearFormats()
    $ws.Range("A4").Value = "Carvão"

    $ws.Range("A5").ClearFormats()
    $ws.Range("A5").Value = "Nuclear"

    $ws.Range("A6").ClearFormats()
    $ws.Range("A6").Value = "Óleos Comb"

    $ws.Range("A7").ClearFormats()
    $ws.Range("A7").Value = "Biomassa"

    $ws.Range("A8").ClearFormats()
    $ws.Range("A8").Value = "Eólica"

    $ws.Range("A9").ClearFormats()
    $ws.Range("A9").Value = "Solar"

    $ws.Range("A10").ClearFormats()
    $ws.Range("A10").Value = "Outros"

    $ws.Range("A11").ClearFormats()
    $ws.Range("A11").Value = "Pot. Compl."

    $ws.Range("A12").ClearFormats()
    $ws.Range("A12").Value = "GD"
}

# ---------------------------------------------------------------------------
# Sheet 5: "Emissoes Totais (MtCO2eq)" -> Atendimento a Ponta style table.
# Add header "Periodo" in A1, fix row labels, drop the "Teto" row entirely.
# ---------------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item(5)

$ws5.Range("B1").Copy()
$ws5.Range("A1").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws5.Range("A1").Value = "Período"

$ws5.Range("A2").ClearFormats()
$ws5.Range("A2").Value = "P.Médio"

$ws5.Range("A3").ClearFormats()
$ws5.Range("A3").Value = "P.Crítico"

$ws5.Range("A4:E4").EntireRow.Delete()

# ---------------------------------------------------------------------------
# Sheet 6: "Custo Total (bilhões de R$)".
# Add header "Tipo Expansão" in A1, change B1 from "Custo" to "2015",
# fix row labels and update the two cost figures.
# ---------------------------------------------------------------------------
$ws6 = $wb.Worksheets.Item(6)

$ws6.Range("B1").Copy()
$ws6.Range("A1").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws6.Range("A1").Value = "Tipo Expansão"

# "Custo" -> "2015" as literal text (matching the other sheets' year
# headers, which are text, not numbers). Round-tripping through a TEXT()
# formula in a scratch cell keeps the value a string without Excel
# reinterpreting "2015" as a number and without leaving a new number
# format behind in the style table.
$scratch = $ws6.Range("Z1")
$scratch.Formula = "=TEXT(2015,""0"")"
$scratch.Copy()
$ws6.Range("B1").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$scratch.Clear()

$ws6.Range("A2").ClearFormats()
$ws6.Range("A2").Value = "Expansão Centralizada"
$ws6.Range("B2").Value = 605

$ws6.Range("A3").ClearFormats()
$ws6.Range("A3").Value = "Expansão por GD"
$ws6.Range("B3").Value = 99
